$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, pushing the current row 57
# (phone "09876543" / "2025-08-12" / 0) down to row 58.
$ws.Rows.Item(57).Insert()

# New row 57: phone 9876543 as a real number, same birthday text,
# 0 points.
$ws.Cells.Item(57, 1).Value = 9876543
$ws.Cells.Item(57, 3).Value = 0

$b57 = $ws.Cells.Item(57, 2)
$b57.NumberFormat = "@"
$b57.Value = "2025-08-12"
$b57.Style = "Normal"

# Former row 57 (now row 58) loses its birthday value; phone stays as
# the zero-padded text "09876543", points stay 0.
$ws.Cells.Item(58, 2).Value = ""
